$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that currently sits in the
# empty paragraph just above "Kinh gui: {donViXacMinh}".
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: "...co ten trong danh sach tot nghiep cac Ky thi: Tu tai, Tot
# nghiep Trung hoc pho thong tai cac hoi dong thi cua {" becomes
# "...co ten trong danh sach tot nghiep thuoc cac Ky thi: {maHeDaoTao} tai
# cac hoi dong thi cua {", and the "_GoBack" bookmark is re-created right in
# the middle of the new "{maHeDaoTao}" placeholder (between "ma" and
# "HeDaoTao}").
# ---------------------------------------------------------------------------
$oldText2 = "uong}/…… (…../…) người (danh sách đính kèm) có tên trong danh sách tốt nghiệp các Kỳ thi: Tú tài, Tốt nghiệp Trung học phổ thông tại các hội đồng thi của {"
$newText2 = "uong}/…… (…../…) người (danh sách đính kèm) có tên trong danh sách tốt nghiệp thuộc các Kỳ thi: {maHeDaoTao} tại các hội đồng thi của {"

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null

$rngMa = $d.Content
$rngMa.Find.ClearFormatting()
$rngMa.Find.Execute("{ma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $rngMa.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Change 3: "... trong danh sach tot nghiep cac Ky thi: Tu tai, Tot nghiep
# Trung hoc pho thong tai cac hoi dong thi cua " becomes "... trong danh
# sach tot nghiep thuoc cac Ky thi: {tenHeDaoTao} tai cac hoi dong thi cua "
# ---------------------------------------------------------------------------
$oldText3 = " trong danh sách tốt nghiệp các Kỳ thi: Tú tài, Tốt nghiệp Trung học phổ thông tại các hội đồng thi của "
$newText3 = " trong danh sách tốt nghiệp thuộc các Kỳ thi: {tenHeDaoTao} tại các hội đồng thi của "

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Execute($oldText3, $true, $false, $false, $false, $false, $true, 1, $false, $newText3, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: "DANH SACH HOC SINH TOT NGHIEP THPT" becomes
# "DANH SACH HOC SINH TOT NGHIEP {maHeDaoTao}" (placeholder bold, like the
# rest of the title).
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$rng4.Find.Execute("DANH SÁCH HỌC SINH TỐT NGHIỆP THPT", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start4 = $rng4.Start
$prefix4 = "DANH SÁCH HỌC SINH TỐT NGHIỆP "
$rng4.Text = $prefix4
$rng4.Font.Bold = $true
$rng4.Font.Color = 0
$rng4.Font.Size = 14
$insPos4 = $start4 + $prefix4.Length
$insRng4 = $d.Range($insPos4, $insPos4)
$insRng4.InsertAfter("{maHeDaoTao}")
$insRng4.Font.Bold = $true
$insRng4.Font.Color = 0
$insRng4.Font.Size = 14

# ---------------------------------------------------------------------------
# Change 5: "DANH SACH TOT NGHIEP THPT TAI CAC HOI DONG THI" becomes
# "DANH SACH TOT NGHIEP {maHeDaoTao}TAI CAC HOI DONG THI".
# ---------------------------------------------------------------------------
$rng5 = $d.Content
$rng5.Find.ClearFormatting()
$rng5.Find.Execute("DANH SÁCH TỐT NGHIỆP THPT TẠI CÁC HỘI ĐỒNG THI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start5 = $rng5.Start
$prefix5 = "DANH SÁCH TỐT NGHIỆP "
$suffix5 = "TẠI CÁC HỘI ĐỒNG THI"
$rng5.Text = $prefix5 + $suffix5
$rng5.Font.Bold = $true
$rng5.Font.Color = 0
$rng5.Font.Size = 14
$insPos5 = $start5 + $prefix5.Length
$insRng5 = $d.Range($insPos5, $insPos5)
$insRng5.InsertAfter("{maHeDaoTao}")
$insRng5.Font.Bold = $true
$insRng5.Font.Color = 0
$insRng5.Font.Size = 14

Write-Host "edit.ps1 completed"
